$d = $word.ActiveDocument

# 1. Remove the "Secundario: Base de Datos de Swapply" paragraph entirely
#    (it directly follows the "Principal: Usuario" paragraph in the Actores cell)
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "*Secundario: Base de Datos de Swapply*") {
        $p.Range.Delete()
        break
    }
}

# 2. Merge the split runs "El usuario no introduce un" + "a" + " " + "descripción " + "para el producto."
#    into a single run of text "El usuario no introduce una descripción para el producto."
$d.Content.Find.Execute(
    "El usuario no introduce una descripción para el producto.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "El usuario no introduce una descripción para el producto.", 2
) | Out-Null

# 3. Merge the split runs "El usuario no " + "sube ninguna foto del producto."
#    into a single run of text "El usuario no sube ninguna foto del producto."
$d.Content.Find.Execute(
    "El usuario no sube ninguna foto del producto.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "El usuario no sube ninguna foto del producto.", 2
) | Out-Null
